$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap the contents of several columns between row 6 and row 7.
# Capture row 6 values first (use Value2 - Value getter has interop quirks).
$A6 = $ws.Range("A6").Value2
$Q6 = $ws.Range("Q6").Value2
$R6 = $ws.Range("R6").Value2
$AC6 = $ws.Range("AC6").Value2
$AM6 = $ws.Range("AM6").Value2
$AO6 = $ws.Range("AO6").Value2

# Capture row 7 values.
$A7 = $ws.Range("A7").Value2
$Q7 = $ws.Range("Q7").Value2
$R7 = $ws.Range("R7").Value2
$AC7 = $ws.Range("AC7").Value2
$AM7 = $ws.Range("AM7").Value2
$AO7 = $ws.Range("AO7").Value2

# Write row 7's original values into row 6.
$ws.Range("A6").Value = $A7
$ws.Range("Q6").Value = $Q7
$ws.Range("R6").Value = $R7
$ws.Range("AC6").Value = $AC7
$ws.Range("AM6").Value = $AM7
$ws.Range("AO6").Value = $AO7

# Write row 6's original values into row 7.
$ws.Range("A7").Value = $A6
$ws.Range("Q7").Value = $Q6
$ws.Range("R7").Value = $R6
$ws.Range("AC7").Value = $AC6
$ws.Range("AM7").Value = $AM6
$ws.Range("AO7").Value = $AO6
